# "little changes in the excel file"
# - append unit suffixes to the suspension-parameter labels in column A
# - widen column A to fit the new (longer) labels
# - move the active selection to B8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$ws.Range("A2").Value = "Sprung mass (ms) [kg]"
$ws.Range("A3").Value = "Unsprung mass (mus) [kg]"
$ws.Range("A4").Value = "stiffness of unsprung (kus) [N/m]"
$ws.Range("A5").Value = "stiffness of sprung (ks) [N/m]"
$ws.Range("A6").Value = "dumping of unsprung (cus) [Ns/m]"
$ws.Range("A7").Value = "dumping of sprung (cs) [Ns/m]"

$ws.Columns.Item(1).ColumnWidth = 30.15

$ws.Range("B8").Select()
